# --------------------------------------------------------------------------
# Adds a new "localdb" command-category column to the hidden '#system' sheet,
# shifting the existing N:AC columns to O:AD, inserts the "localdb" entry
# into the "target" category list (column A), and fixes up all the
# definedNames so that they keep pointing at the right ranges.
# --------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

$oldLastCol = 29   # AC
$newLastCol = 30   # AD
$lastRow    = 127

# 1. Read the whole existing grid (A1:AC127) in one shot.
$src = $ws.Range("A1:AC127").Value2

# 2. Build the new grid (A1:AD127): columns 1..13 (A..M) stay where they
#    are, column 14 (N) becomes the brand-new "localdb" column, and the
#    old columns 14..29 (N..AC) move right by one slot to 15..30 (O..AD).
$dst = New-Object 'object[,]' $lastRow,$newLastCol

for ($r = 1; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le 13; $c++) {
        $dst[$r-1, $c-1] = $src[$r, $c]
    }
    for ($c = 14; $c -le $oldLastCol; $c++) {
        $dst[$r-1, $c] = $src[$r, $c]
    }
}

# 3. Populate the new "localdb" column (column 14 / "N").
$localdb = @(
    "localdb",
    "cloneTable(var,source,target)",
    "dropTables(var,tables)",
    "exportCSV(sql,output)",
    "importRecords(var,sourceDb,sql,table)",
    "purge(var)",
    "runSQLs(var,sqls)"
)
for ($i = 0; $i -lt $localdb.Length; $i++) {
    $dst[$i, 13] = $localdb[$i]
}

# 4. Insert "localdb" into the "target" list held in column A: rows 14-29
#    shift down to 15-30, and row 14 becomes "localdb".
for ($r = 29; $r -ge 14; $r--) {
    $dst[$r, 0] = $dst[$r-1, 0]
}
$dst[13, 0] = "localdb"

# 5. Write the whole new grid back out.
$ws.Range("A1:AD127").Value = $dst

# 6. Fix up the definedNames that pointed into the shifted columns (and the
#    "target" range, which now runs one row further), then add the new
#    "localdb" name.
function Set-SystemName($name, $colLetter, $firstRow, $lastRowForName) {
    $wb.Names.Item($name).Delete()
    $refersTo = "='#system'!`$$colLetter`$$firstRow`:`$$colLetter`$$lastRowForName"
    $wb.Names.Add($name, $refersTo)
}

Set-SystemName "macro"      "O"  2 4
Set-SystemName "mail"       "P"  2 2
Set-SystemName "number"     "Q"  2 16
Set-SystemName "pdf"        "R"  2 16
Set-SystemName "rdbms"      "S"  2 7
Set-SystemName "redis"      "T"  2 10
Set-SystemName "sms"        "U"  2 2
Set-SystemName "sound"      "V"  2 5
Set-SystemName "ssh"        "W"  2 9
Set-SystemName "step"       "X"  2 4
Set-SystemName "target"     "A"  2 30
Set-SystemName "web"        "Y"  2 127
Set-SystemName "webalert"   "Z"  2 8
Set-SystemName "webcookie"  "AA" 2 8
Set-SystemName "ws"         "AB" 2 17
Set-SystemName "ws.async"   "AC" 2 8
Set-SystemName "xml"        "AD" 2 21

$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")
